$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.640.59"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.650.09"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.08"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.16"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +6.88%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("E11").Value = "  +6.66%  "
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.127.51"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.613.41"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.675.95"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.99%  "
$ws.Range("E18").Value = "  +3.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.79"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.48"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +3.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.72"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +7.15%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.992"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +4.99%  "
$ws.Range("E29").Value = "  +8.00%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("E32").Value = "  +5.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.32"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.30"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.11"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +6.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.917"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +7.36%  "
$ws.Range("E37").Value = "  +8.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.914"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +13.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.68"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  +8.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "309.35"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +10.68%  "
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0982"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +5.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0551"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.50"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.06"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +12.50%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.66"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  +5.49%  "
$ws.Range("E51").Value = "  +7.78%  "
